# story number: 10,21 - cancel sale, shir bar , 31.12.18
#
# The "size" column (C: product code | name | size | amount | color | price)
# is removed entirely, so amount/color/price shift one column to the left
# (account | color | price end up in C/D/E). The "amount" header is renamed
# to "account", the trailing space in the "price " header is dropped, a
# couple of amount/color values are updated, and the two "shirt" rows become
# "shirts".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the whole "size" column; amount/color/price shift left into C/D/E
# and the now-unused "size" / "s/m/l" / "s/m" shared strings drop out.
$ws.Columns("C").Delete()

# Header renames: "amount" -> "account", "price " -> "price".
$ws.Range("C1").Value = "account"
$ws.Range("E1").Value = "price"

# Updated amounts for shirts (row 2) and jeans (row 3).
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 14

# "shirt" -> "shirts" for the last two products.
$ws.Range("B6").Value = "shirts"
$ws.Range("B7").Value = "shirts"

# Color updates: jeans -> blue, coats -> brown, last shirts row -> white.
$ws.Range("D3").Value = "blue"
$ws.Range("D5").Value = "brown "
$ws.Range("D7").Value = "white"
